# Insert two new "Farewell" entries ("Have a nice day" / "Bye") right after
# the existing Farewell block (row 110), before the "Options" block which
# currently starts at row 111. This shifts the old rows 111-121 down to
# 113-123.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 111, pushing existing rows 111-121 down to 113-123
$ws.Range("A111:B112").EntireRow.Insert()

# Fill the two new rows with the new Farewell sentences
$ws.Range("A111").Value = "Farewell"
$ws.Range("B111").Value = "Have a nice day"
$ws.Range("A112").Value = "Farewell"
$ws.Range("B112").Value = "Bye"

# Update the view to match the saved state (active cell / scroll position / zoom)
$ws.Application.ActiveWindow.Zoom = 160
$ws.Application.ActiveWindow.ScrollRow = 94
$ws.Range("B113").Select()
